$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.325.19"
$ws.Range("E2").Value = "  -3.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.71"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.95"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4622"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3942"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.91"
$ws.Range("E9").Value = "  -11.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07943"
$ws.Range("E10").Value = "  -5.72%  "
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.47"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.922"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.821.10"
$ws.Range("E14").Value = "  -7.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.129"
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.12"
$ws.Range("E17").Value = "  -3.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001030"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06544"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.470"
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.330.41"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.069.73"
$ws.Range("E26").Value = "  -5.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.44"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.61"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.063"
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.465"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.93"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09445"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9472"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.434"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.585"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.254"
$ws.Range("E36").Value = "  -5.12%  "
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02227"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.212"
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.022"
$ws.Range("E41").Value = "  -9.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5932"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.18"
$ws.Range("E44").Value = "  -7.73%  "
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5618"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.427"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.920"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06760"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.37"
$ws.Range("E51").Value = "  -1.19%  "
